$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices / volume%) per latest GitHub Actions refresh

$ws.Range("D2").Value = "39.596.21"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "2.161.56"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.79"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0844"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.88"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "2.480.68"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "2.166.62"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "39.563.75"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.81"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.84%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.95%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("B39").Value = "BinanceUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("B40").Value = "FTXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.96"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +19.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.53"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0227"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").Value = "1.514.27"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +28.76%  "
$ws.Range("E51").Value = "  +0.46%  "
